# Update profit files after running on 2025-10-17
# Appends the newly computed allocation row (Date, BTC, KAS) for 10/17/2025
# at the bottom of the existing table on Sheet1 (row 46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-like string to be stored as literal
# text (matching the existing "Date" column cells, which are plain text,
# not real Excel dates) instead of being auto-converted to a date serial
# number. Resetting the style back to "Normal" afterwards clears the
# quote-prefix formatting Excel applies when it detects the forced-text
# input, so the new cell ends up unstyled just like its neighbours.
$ws.Range("A46").Value = "'10/17/2025"
$ws.Range("A46").Style = "Normal"

$ws.Range("B46").Value = 0.1864789165136171
$ws.Range("C46").Value = 0.8135210834863829
